$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: update the title text
$ws.Range("A1").Value = "Sample Excel Data test 2"

# A2: convert the numeric value 12345 into the text value "12345".
# Temporarily mark the cell as Text so Excel doesn't re-coerce the
# numeric-looking string back into a number, then restore the
# cell's (unformatted) Normal style so no stray number format sticks.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "12345"
$ws.Range("A2").Style = "Normal"
